$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5 (shifts existing rows 5.. down by one)
$ws.Rows(5).Insert()

# Populate the newly inserted row with the new daily record
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 24909.02
$ws.Cells.Item(5, 3).Value = 8
$ws.Cells.Item(5, 4).Value = 2025
$ws.Cells.Item(5, 5).Value = "08/2025"
